$d = $word.ActiveDocument
$r = $d.Content
$found = $r.Find.Execute("https://lidorprototype.github.io/", $true, $false, $false, $false, $false,
                         $true, 1, $false, "", 0)
$r.Collapse(0)
$newText = [char]9 + "(you can find everything in here)"
$r.InsertAfter($newText)

$r2 = $d.Content
$searchText = [char]9 + "(you can find everything in here)"
$found2 = $r2.Find.Execute($searchText, $true, $false, $false, $false, $false,
                         $true, 1, $false, "", 0)

$r2.Style = "Hyperlink"
$r2.Font.NameAscii = "Cambria Math"
$r2.Font.Name = "Cambria Math"
$r2.Font.NameFarEast = "Adobe Ming Std L"
$r2.Font.Size = 10
$r2.Font.SizeBi = 10
$r2.Font.Underline = 0
$r2.Font.Color = 0

Write-Output "done"
